# Update CMPA template so that all sample dates fall within the same month
# (July 2023) and refresh the sample internet-media domain from
# "topky.sk" to "zoznam.sk" (rows 12-13, column B) on the "Spoty" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spoty")

# --- Date columns (F = start date, G = end date) -----------------------
# Row 3-11: single date column F
$ws.Cells.Item(3, 6).Value  = 45127   # F3  2023-07-20
$ws.Cells.Item(4, 6).Value  = 45128   # F4  2023-07-21
$ws.Cells.Item(5, 6).Value  = 45129   # F5  2023-07-22
$ws.Cells.Item(6, 6).Value  = 45127   # F6  2023-07-20
$ws.Cells.Item(7, 6).Value  = 45128   # F7  2023-07-21
$ws.Cells.Item(8, 6).Value  = 45130   # F8  2023-07-23
$ws.Cells.Item(9, 6).Value  = 45126   # F9  2023-07-19
$ws.Cells.Item(10, 6).Value = 45127   # F10 2023-07-20
$ws.Cells.Item(11, 6).Value = 45129   # F11 2023-07-22

# Rows 12-18: start date (F) + end date (G)
$ws.Cells.Item(12, 6).Value = 45119   # F12 2023-07-12
$ws.Cells.Item(12, 7).Value = 45127   # G12 2023-07-20

$ws.Cells.Item(13, 6).Value = 45110   # F13 2023-07-03
$ws.Cells.Item(13, 7).Value = 45110   # G13 2023-07-03

$ws.Cells.Item(14, 6).Value = 45127   # F14 2023-07-20
$ws.Cells.Item(14, 7).Value = 45129   # G14 2023-07-22

$ws.Cells.Item(15, 6).Value = 45116   # F15 2023-07-09
$ws.Cells.Item(15, 7).Value = 45119   # G15 2023-07-12

$ws.Cells.Item(16, 6).Value = 45116   # F16 2023-07-09
$ws.Cells.Item(16, 7).Value = 45119   # G16 2023-07-12

$ws.Cells.Item(17, 6).Value = 45116   # F17 2023-07-09
$ws.Cells.Item(17, 7).Value = 45119   # G17 2023-07-12

$ws.Cells.Item(18, 6).Value = 45116   # F18 2023-07-09
$ws.Cells.Item(18, 7).Value = 45119   # G18 2023-07-12

# --- Sample domain swap: topky.sk -> zoznam.sk (rows 12-13, column B) --
$ws.Cells.Item(12, 2).Value = "zoznam.sk"
$ws.Cells.Item(13, 2).Value = "zoznam.sk"
